$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'21.330.19"
$ws.Range("E2").Value = "'  +4.19%  "
$ws.Range("D3").Value = "'1.546.63"
$ws.Range("E3").Value = "'  +5.07%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("D5").Value = "'0.9701"
$ws.Range("E5").Value = "'  +0.78%  "
$ws.Range("D6").Value = "'282.19"
$ws.Range("E6").Value = "'  +2.28%  "
$ws.Range("D7").Value = "'0.3630"
$ws.Range("E7").Value = "'  -0.58%  "
$ws.Range("D8").Value = "'0.3197"
$ws.Range("E8").Value = "'  +4.42%  "
$ws.Range("D9").Value = "'40.98"
$ws.Range("E9").Value = "'  +2.85%  "
$ws.Range("D10").Value = "'1.112"
$ws.Range("E10").Value = "'  +6.22%  "
$ws.Range("D11").Value = "'0.06916"
$ws.Range("E11").Value = "'  +4.67%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "'  -0.07%  "
$ws.Range("D13").Value = "'5.719"
$ws.Range("E13").Value = "'  +4.91%  "
$ws.Range("D14").Value = "'18.94"
$ws.Range("E14").Value = "'  +4.00%  "
$ws.Range("D15").Value = "'6.413"
$ws.Range("E15").Value = "'  +4.08%  "
$ws.Range("D16").Value = "'0.00001054"
$ws.Range("E16").Value = "'  +2.48%  "
$ws.Range("D17").Value = "'0.9695"
$ws.Range("E17").Value = "'  -0.18%  "
$ws.Range("D18").Value = "'1.546.60"
$ws.Range("E18").Value = "'  +4.95%  "
$ws.Range("D19").Value = "'0.06152"
$ws.Range("E19").Value = "'  +4.50%  "
$ws.Range("D20").Value = "'73.09"
$ws.Range("E20").Value = "'  +5.91%  "
$ws.Range("D21").Value = "'5.744"
$ws.Range("E21").Value = "'  +5.42%  "
$ws.Range("D22").Value = "'15.27"
$ws.Range("E22").Value = "'  +6.01%  "
$ws.Range("D23").Value = "'11.42"
$ws.Range("D24").Value = "'2.323"
$ws.Range("E24").Value = "'  +3.31%  "
$ws.Range("D25").Value = "'21.375.67"
$ws.Range("E25").Value = "'  +4.15%  "
$ws.Range("D26").Value = "'147.82"
$ws.Range("E26").Value = "'  +4.25%  "
$ws.Range("D27").Value = "'2.284"
$ws.Range("E27").Value = "'  +7.37%  "
$ws.Range("E28").Value = "'  +3.72%  "
$ws.Range("D29").Value = "'1.717.84"
$ws.Range("E29").Value = "'  +5.43%  "
$ws.Range("D30").Value = "'118.61"
$ws.Range("E30").Value = "'  +4.65%  "
$ws.Range("D31").Value = "'4.053"
$ws.Range("E31").Value = "'  +4.42%  "
$ws.Range("D32").Value = "'0.8743"
$ws.Range("E32").Value = "'  +8.38%  "
$ws.Range("D33").Value = "'5.256"
$ws.Range("E33").Value = "'  +6.19%  "
$ws.Range("D34").Value = "'0.08053"
$ws.Range("E34").Value = "'  +2.36%  "
$ws.Range("D35").Value = "'1.539"
$ws.Range("E35").Value = "'  +0.73%  "
$ws.Range("E36").Value = "'  +5.41%  "
$ws.Range("D37").Value = "'1.218"
$ws.Range("E37").Value = "'  -3.33%  "
$ws.Range("D38").Value = "'0.05875"
$ws.Range("E38").Value = "'  +2.55%  "
$ws.Range("D39").Value = "'0.2008"
$ws.Range("E39").Value = "'  +7.02%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.02126"
$ws.Range("E40").Value = "'  +4.31%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.964"
$ws.Range("E41").Value = "'  +4.19%  "
$ws.Range("D42").Value = "'10.78"
$ws.Range("E42").Value = "'  +3.67%  "
$ws.Range("D43").Value = "'0.9694"
$ws.Range("E43").Value = "'  +0.46%  "
$ws.Range("D44").Value = "'0.5512"
$ws.Range("E44").Value = "'  +4.52%  "
$ws.Range("D45").Value = "'12.60"
$ws.Range("E45").Value = "'  +4.87%  "
$ws.Range("D46").Value = "'3.581"
$ws.Range("E46").Value = "'  +2.26%  "
$ws.Range("D47").Value = "'0.5503"
$ws.Range("E47").Value = "'  +6.50%  "
$ws.Range("D48").Value = "'122.39"
$ws.Range("E48").Value = "'  +4.84%  "
$ws.Range("D49").Value = "'1.880"
$ws.Range("E49").Value = "'  +6.52%  "
$ws.Range("D50").Value = "'0.06611"
$ws.Range("E50").Value = "'  +2.62%  "
$ws.Range("D51").Value = "'70.13"
$ws.Range("E51").Value = "'  +5.00%  "
